$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'44.136.06"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.96%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.246.71"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.43%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.11%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'268.03"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +4.30%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'87.56"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +12.47%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.620"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.46%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.05%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.615"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +3.39%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'45.89"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +6.39%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +2.21%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = "'  +8.20%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  +2.47%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'2.584.08"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.30%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'15.01"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +3.94%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'2.247.57"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.39%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.800"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.01%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'44.100.73"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +3.00%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.59%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'6.05"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.17%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'70.33"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  -1.05%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'2.40"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +5.03%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'234.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.76%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'8.95"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -3.63%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('B25').Value = "'Dai"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'1.00"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.03%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('B26').Value = "'PancakeSwap"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'2.55"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +16.12%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'11.01"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +2.68%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'3.58"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +7.09%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'40.66"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -4.85%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'2.27"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +2.53%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'175.54"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.07%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.0915"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +4.81%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.92%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = "'  +4.22%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D36').Value = "'0.112"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +3.86%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.64%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'4.40"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.28%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +16.99%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'12.78"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -2.83%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'2.15"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +1.71%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'65.40"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +6.50%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.56%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'5.41"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.82%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.0992"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +2.05%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  -0.68%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'100.46"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -2.81%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'1.21"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +6.84%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +1.83%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').Value = "'  -9.28%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'1.52"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +3.36%  "
$ws.Range('E51').Style = 'Normal'
